$wb = $excel.ActiveWorkbook

# --- "Folder Inventory" sheet: a new, more-recently-updated folder
# ("Microsoft Azure AI Agents") bubbled to the top of the list (row 3,
# just below the header + the "Azure Landing Zone" row), pushing every
# other row from 3..36 down by one. Rows 37+ are untouched.
$ws = $wb.Worksheets.Item("Folder Inventory")

# Shift the old rows 3..35 down into rows 4..36 (Folder Path / Folder
# Name / Last Updated columns only - File Count & Parent Directory are
# identical for every row so they don't need touching).
$shifted = $ws.Range("A3:C35").Value2
$ws.Range("A4:C36").Value2 = $shifted

# Write the new top entry into row 3.
$ws.Range("A3").Value2 = "Microsoft Azure AI Agents"
$ws.Range("B3").Value2 = "Microsoft Azure AI Agents"
$ws.Range("C3").Value2 = "2025-06-11 20:13:48 +0530"

# --- "Metadata" sheet: refresh the generation timestamp + workflow run
# counter.
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value2 = "2025-06-11 15:41:19 UTC"

# "Workflow Run" (B5) is stored as text in the original workbook (unlike
# "Total Folders" in B4, which is a real number), so force text entry to
# avoid Excel auto-converting "4" into a numeric cell.
$meta.Range("B5").NumberFormat = "@"
$meta.Range("B5").Value2 = "4"
$meta.Range("B5").ClearFormats()
